# "lagt til nye forslag på den nye tidstabel (excel)"
# Adds new route/price suggestions to the timetable and switches the
# "klokkeslett" (time) column from a fixed h:mm time format to plain
# General formatting (some rows now hold free-text time strings instead
# of real Excel time serials).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a piece of text into a cell WITHOUT letting Excel's
# automatic "number-looking text becomes a real number" conversion
# kick in (used for strings like "10.00" that must stay text). We stage
# the literal value in a scratch cell that is explicitly formatted as
# Text, copy it, and paste-special just the values into the real
# target - this carries over the String type but none of the scratch
# cell's formatting (so the destination is left with its own/default
# style, matching a plain text entry).
# ---------------------------------------------------------------------
function Set-SafeText {
    param($range, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# klokkeslett column: drop the dedicated time-number format (numFmtId
# 20, "h:mm") in favour of General on every data row that uses it.
# ---------------------------------------------------------------------
$ws.Range("D2:D15").NumberFormat = "General"

# D6 used to hold the text "08:00, 17:00"; it is now a real time serial.
$ws.Range("D6").Value = 0.33353206018518522

# ---------------------------------------------------------------------
# New row: Sandefjord -> Strømstad / Strømstad -> Sandefjord
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Sandefjord"
$ws.Range("B14").Value = "Strømstad"
$ws.Range("C14").Value = "mandag, onsdag, fredag"
Set-SafeText $ws.Range("D14") "8.00, 17.00"
$ws.Range("E14").Value = "NOK 449.00"
$ws.Range("F14").Value = "NOK 299.00"
$ws.Range("G14").Value = "NOK 100.00"
$ws.Range("H14").Value = "NOK 99.00"

$ws.Range("A15").Value = "Strømstad"
$ws.Range("B15").Value = "Sandefjord"
Set-SafeText $ws.Range("D15") "9.00, 18.00"

# ---------------------------------------------------------------------
# New row: Stavanger -> Bergen / Bergen -> Stavanger
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Stavanger"
$ws.Range("B17").Value = "Bergen"
$ws.Range("C17").Value = "tirsdag, torsdag"
Set-SafeText $ws.Range("D17") "8.00, 15.00"
$ws.Range("E17").Value = "NOK 549.00"
$ws.Range("F17").Value = "NOK 399.00"
$ws.Range("G17").Value = "NOK 100.00"
$ws.Range("H17").Value = "NOK 149.00"

$ws.Range("A18").Value = "Bergen"
$ws.Range("B18").Value = "Stavanger"
Set-SafeText $ws.Range("D18") "8.00, 16.00"

# ---------------------------------------------------------------------
# New row: Oslo -> Kiel / Kiel -> Oslo
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Oslo"
$ws.Range("B20").Value = "Kiel"
$ws.Range("C20").Value = "lørdag"
Set-SafeText $ws.Range("D20") "10.00"
$ws.Range("E20").Value = "NOK 749.00"
$ws.Range("F20").Value = "NOK 549.00"
$ws.Range("G20").Value = "NOK 100.00"
$ws.Range("H20").Value = "NOK 200.00"

$ws.Range("A21").Value = "Kiel"
$ws.Range("B21").Value = "Oslo"
Set-SafeText $ws.Range("D21") "10.00"
